$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("B4,B12")
$rng.Interior.ThemeColor = 10
$rng.Interior.TintAndShade = 0.39997558519241921
